# "update ui when switch cur player"
# CardPile.xlsx: drop the old multi-run emoji note from C1, and add a new
# "ResName" column (E) holding the pinyin id for each basic/common card's
# Chinese Name (column B). Rows beyond the basic card catalogue (weapons,
# armor, horses, generic equipment ids >= 133) are left without a ResName,
# matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 used to hold a single rich-text note cell (C1) explaining the
# Suit encoding (1:♦ 2:♥ 3:♣ 4:♠). That note is removed entirely so row 1
# goes back to being fully blank. Deleting then re-inserting the row drops
# every trace of it (height, cell) while leaving all other row numbers
# untouched.
$ws.Rows("1:1").Delete()
$ws.Rows("1:1").Insert()

# --- New column E: "ResName" header + "string" dtype row, mirroring the
# existing Id/Name/Suit/Rank header (row 2) and dtype (row 3) layout.
$ws.Range("E2").Value = "ResName"
$ws.Range("E3").Value = "string"

# --- ResName values: the pinyin transliteration of each card's Chinese
# Name (column B), filled contiguously per card-name block.
$ws.Range("E4:E47").Value = "sha"
$ws.Range("E48:E71").Value = "shan"
$ws.Range("E72:E83").Value = "tao"
$ws.Range("E84:E88").Value = "jiu"
$ws.Range("E89").Value = "taoyuan"
$ws.Range("E90").Value = "wanjian"
$ws.Range("E91:E92").Value = "jiedao"
$ws.Range("E93:E94").Value = "wugu"
$ws.Range("E95:E97").Value = "nanman"
$ws.Range("E98:E100").Value = "juedou"
$ws.Range("E101:E103").Value = "huogong"
$ws.Range("E104:E107").Value = "wuzhong"
$ws.Range("E108:E112").Value = "shunshou"
$ws.Range("E113:E118").Value = "guochai"
$ws.Range("E119:E124").Value = "tiesuo"
$ws.Range("E125:E131").Value = "wuxie"
$ws.Range("E132:E133").Value = "shandian"
$ws.Range("E134:E135").Value = "bingliang"

# The filled-in ResName cells picked up a (no-op) alignment style along the
# way -- reproduce that so the style table shape matches.
$ws.Range("E4:E135").WrapText = $false

# --- Column widths: with the wide note gone, columns A-C re-fit to their
# (now much narrower) content; column D was already sized and is untouched.
$ws.Columns("A:C").AutoFit()

# --- Leave the view roughly where editing finished.
$ws.Activate()
$ws.Range("H139").Select()
$excel.ActiveWindow.ScrollRow = 128
$excel.ActiveWindow.ScrollColumn = 1
